# Rename title in ppt: "Indeksi kod PostgreSQL" -> "Indeksi u PostgreSQL"
#
# The original title run reads "Indeksi kod PostgreSQL". The word "kod" (plus
# the trailing space that follows it) is replaced with "u " so the title
# becomes "Indeksi u PostgreSQL", split across three runs exactly the way
# PowerPoint splits a run when you retype a word in the middle of existing
# text: "Indeksi " | "u " | "PostgreSQL".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# "Indeksi kod PostgreSQL"
#          ^^^^            -> characters 9..12 are "kod "
$mid = $tr.Characters(9, 4)
$mid.Text = "u "
